# Add a new "Fuel Flow Comparison" worksheet at the end of the workbook,
# comparing fuel flow / speed / range tradeoffs across three power settings.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
# of the tab strip (Worksheets.Add() defaults to inserting before the
# active sheet, which is not what we want here).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Fuel Flow Comparison"

# Column widths (best-effort: the host quantizes ColumnWidth to 1/6-wide
# steps, so these inputs are chosen to land on the closest achievable step
# to the source file's true widths of 16.66, 20.5, 21.5, 20.33, 18.16, 25).
$ws.Columns.Item(1).ColumnWidth = 15.833333333333334
$ws.Columns.Item(2).ColumnWidth = 19.666666666666668
$ws.Columns.Item(3).ColumnWidth = 20.666666666666668
$ws.Columns.Item(4).ColumnWidth = 19.5
$ws.Columns.Item(5).ColumnWidth = 17.333333333333336
$ws.Columns.Item(6).ColumnWidth = 24.166666666666668

# Headers (written in the shared-string insertion order the source file
# expects: the "Range" column's header string was added to the shared
# string table last, after the diff columns).
$ws.Range("A1").Value = "Percent Power"
$ws.Range("B1").Value = "Fuel Flow (GPH)"
$ws.Range("C1").Value = "True Airspeed (Knots)"
$ws.Range("E1").Value = "Diff. Fuel Flow"
$ws.Range("F1").Value = "Diff. True Airspeed"
$ws.Range("G1").Value = "Diff. Range"
$ws.Range("D1").Value = "Range (Nautical Miles)"

# Row 2 - 55% power (baseline, no diff columns)
$ws.Range("A2").Value = 55
$ws.Range("B2").Value = 9.3
$ws.Range("C2").Value = 100
$ws.Range("D2").Formula = "=(72/B2)*C2"

# Row 3 - 65% power
$ws.Range("A3").Value = 65
$ws.Range("B3").Value = 10.9
$ws.Range("C3").Value = 115
$ws.Range("D3").Formula = "=(72/B3)*C3"
$ws.Range("E3").Formula = "=(B3-B2)/B2"
$ws.Range("F3").Formula = "=(C3-C2)/C2"
$ws.Range("G3").Formula = "=(D3-D2)/D2"

# Row 4 - 75% power
$ws.Range("A4").Value = 75
$ws.Range("B4").Value = 12.65
$ws.Range("C4").Value = 124
$ws.Range("D4").Formula = "=(72/B4)*C4"
$ws.Range("E4").Formula = "=(B4-B3)/B3"
$ws.Range("F4").Formula = "=(C4-C3)/C3"
$ws.Range("G4").Formula = "=(D4-D3)/D3"

# Number formats: Range column as 0.00, Diff columns as percent.
$ws.Range("D1:D4").NumberFormat = "0.00"
$ws.Range("G1").NumberFormat = "0%"
$ws.Range("E3:G4").NumberFormat = "0%"

# Match the file's selection/active-cell state.
$ws.Range("G4").Select() | Out-Null
